$wb = $excel.ActiveWorkbook

# --- Populate Sheet2 with the BatchMode list ---
$ws2 = $wb.Worksheets.Item("Sheet2")
$ws2.Range("A1").Value = "BatchMode"
$ws2.Range("A2").Value = "PMKVY-RPL"
$ws2.Range("A3").Value = "PMKVY-STT"
$ws2.Range("A4").Value = "PMKVY-Special"
$ws2.Range("A5").Value = "NSKFDC"
$ws2.Range("A6").Value = "NBCFDC"
$ws2.Range("A7").Value = "CB_Scheme"
$ws2.Range("A8").Value = "State Skill"
$ws2.Range("A9").Value = "MNRE"
$ws2.Range("A10").Value = "MM"

$ws2.Range("A1").HorizontalAlignment = -4108
$ws2.Range("A1").Font.Bold = $true
$ws2.Columns.Item(1).ColumnWidth = 14

$ws2.Range("B8").Select()

# --- Define the named range "BatchMode" referring to the list ---
$wb.Names.Add("BatchMode", "=Sheet2!`$A`$2:`$A`$10")

# --- Sheet1: update H2 text and add data validation to column H ---
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws1.Range("H2").Value = "Select from DropDown List"

$ws1.Range("H3:H1048576").Validation.Add(3, 1, 1, "=BatchMode")
$ws1.Range("H3:H1048576").Validation.InputMessage = ""
$ws1.Range("H3:H1048576").Validation.ErrorMessage = ""
$ws1.Range("H3:H1048576").Validation.ShowInput = $true
$ws1.Range("H3:H1048576").Validation.ShowError = $true

$ws1.Range("H3").Select()
$ws1.Application.ActiveWindow.ScrollColumn = 3
